$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New X / Y data values (rows 2-16, columns A and B)
$data = @(
    @(20,   88.6),
    @(16,   71.6),
    @(19.8, 93.3),
    @(18.4, 84.3),
    @(17.1, 80.6),
    @(15.5, 75.2),
    @(14.7, 69.7),
    @(17.1, 82),
    @(15.4, 69.4),
    @(16.2, 83.3),
    @(15,   79.6),
    @(17.2, 82.6),
    @(16,   80.6),
    @(17,   83.5),
    @(14.4, 76.3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Selection moves from D11 to B1
$ws.Range("B1").Select() | Out-Null
